# Add a new "General Settings" parameter row to the project parameters sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "General Settings"
$ws.Range("B56").Value = "Default IGU Service Lifetime (years)"
$ws.Range("C56").Value = 25
$ws.Range("D56").Value = "Default age assumption for IGUs at end-of-life"
